$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) for every data row (2-31) was stamped with the
# mis-formatted text "5-6-2012-13". Per the commit message, the NBA stats
# for this game were off by one day because of how the stats were
# originally reported, so the column is corrected to the proper
# ISO-style date string "2013-05-06".
#
# Cells are plain text (not real Excel dates), so we build the replacement
# text via a temporary formula and convert it back to a literal value with
# Copy/PasteSpecial(xlPasteValues). This keeps Excel's automatic
# date-recognition (which would otherwise turn "2013-05-06" typed directly
# into a date serial number) from kicking in, while leaving cell
# formatting/styles untouched.
$dateRange = $ws.Range("BF2:BF31")
$dateRange.Formula = '="2013-05-06"'
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
